$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where columns AB:AF (previously all zeros) become blank
$clearRanges = @(
    "AB57:AF57",
    "AB58:AF58",
    "AB71:AF71",
    "AB72:AF72",
    "AB73:AF73",
    "AB77:AF77",
    "AB78:AF78"
)
foreach ($r in $clearRanges) {
    $ws.Range($r).ClearContents()
}

# Rows 64 and 79: everything except column D (already blank) and the tail
# columns (AC:AF on row 64, AC:AF on row 79) becomes blank
$ws.Range("B64:C64").ClearContents()
$ws.Range("E64:AB64").ClearContents()

$ws.Range("B79:C79").ClearContents()
$ws.Range("E79:AB79").ClearContents()

# Small numeric corrections (float rounding updates)
$ws.Range("AF68").Value = 155224.992
$ws.Range("AB70").Value = -96891.016
$ws.Range("X74").Value = -52076.008
$ws.Range("AB74").Value = -86303.992
$ws.Range("X80").Value = -49567
$ws.Range("AB80").Value = -80592.984
